# Add a new "Type" column (I) to the Arraignments sheet, classifying each
# case's charge as Moving / Non-moving / Criminal, for the cost calculator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("I1").Value = "Type"

# Row-by-row classification
$ws.Range("I2").Value  = "Moving"
$ws.Range("I3").Value  = "Non-moving"
$ws.Range("I4").Value  = "Criminal"
$ws.Range("I5").Value  = "Criminal"
$ws.Range("I6").Value  = "Moving"
$ws.Range("I7").Value  = "Moving"
$ws.Range("I8").Value  = "Non-moving"
$ws.Range("I9").Value  = "Moving"
$ws.Range("I10").Value = "Criminal"
$ws.Range("I11").Value = "Moving"
$ws.Range("I12").Value = "Moving"
$ws.Range("I13").Value = "Criminal"
$ws.Range("I14").Value = "Moving"
$ws.Range("I15").Value = "Moving"
$ws.Range("I16").Value = "Moving"
$ws.Range("I17").Value = "Moving"
$ws.Range("I18").Value = "Moving"
$ws.Range("I19").Value = "Non-moving"
$ws.Range("I20").Value = "Non-moving"
$ws.Range("I21").Value = "Moving"
$ws.Range("I22").Value = "Moving"
$ws.Range("I23").Value = "Moving"

# Match the author's final on-screen selection after the edit.
$ws.Range("I24").Select()
